$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes the existing rows 10..35 down to 11..36,
# keeping column widths / row formats consistent with the rest of the data table)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price observation
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44487
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 110
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("N10").Value = "`$/saco 25 kilos"
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 400
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
